$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

Set-TextValue "D2" "27.083.62"
$ws.Range("E2").Value = "  +0.66%  "
Set-TextValue "D3" "1.567.82"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").Value = "  +0.72%  "
Set-TextValue "D5" "208.74"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("E7").Value = "  +0.71%  "
Set-TextValue "D8" "22.03"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("E11").Value = "  +0.64%  "
Set-TextValue "D12" "1.584.26"
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("E14").Value = "  +0.26%  "
Set-TextValue "D15" "27.071.08"
$ws.Range("E15").Value = "  +0.63%  "
Set-TextValue "D16" "61.97"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("E17").Value = "  +0.29%  "
Set-TextValue "D19" "215.75"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("E21").Value = "  +2.18%  "
Set-TextValue "D22" "9.19"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("E23").Value = "  -0.02%  "
Set-TextValue "D24" "154.19"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  -0.32%  "
Set-TextValue "D26" "15.06"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("E29").Value = "  +4.53%  "
$ws.Range("E30").Value = "  +0.93%  "
Set-TextValue "D31" "3.23"
$ws.Range("E31").Value = "  +0.41%  "
Set-TextValue "D32" "3.18"
$ws.Range("E32").Value = "  +2.40%  "
Set-TextValue "D33" "1.427.74"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  +13.02%  "
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("E37").Value = "  +1.09%  "
Set-TextValue "D38" "0.532"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("E41").Value = "  +4.16%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("E43").Value = "  +0.65%  "
Set-TextValue "D44" "64.62"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("E45").Value = "  +0.31%  "
Set-TextValue "D46" "1.704.37"
$ws.Range("E46").Value = "  +1.15%  "
Set-TextValue "D47" "86.63"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("E51").Value = "  +0.63%  "